$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 399.7586
$ws.Range("I28").Value = 439.83334
$ws.Range("J28").Value = 207.4
$ws.Range("K28").Value = 439.83334
$ws.Range("L28").Value = 207.4
$ws.Range("M28").Value = 45.16665999999998
$ws.Range("N28").Value = -1177.4
$ws.Range("H74").Value = 5887368
$ws.Range("I74").Value = 9094772
$ws.Range("J74").Value = 7126.5
$ws.Range("K74").Value = 9094772
$ws.Range("L74").Value = 7126.5
$ws.Range("M74").Value = -9093836
$ws.Range("N74").Value = -8998.5
$ws.Range("H77").Value = 5887368
$ws.Range("I77").Value = 9094772
$ws.Range("J77").Value = 7126.5
$ws.Range("K77").Value = 45473860
$ws.Range("L77").Value = 35632.5
$ws.Range("M77").Value = -45469180
$ws.Range("N77").Value = -44992.5
$ws.Range("H98").Value = 2454.3235
$ws.Range("I98").Value = 1266.6786
$ws.Range("J98").Value = 7996.6665
$ws.Range("K98").Value = 1266.6786
$ws.Range("L98").Value = 7996.6665
$ws.Range("M98").Value = 231.3214
$ws.Range("N98").Value = -10992.6665
$ws.Range("H112").Value = 10001551
$ws.Range("I112").Value = 500000350
$ws.Range("J112").Value = 1575.9183
$ws.Range("K112").Value = 1500001050
$ws.Range("L112").Value = 4727.7549
$ws.Range("M112").Value = -1499999942
$ws.Range("N112").Value = -6943.7549
$ws.Range("H122").Value = 2454.3235
$ws.Range("I122").Value = 1266.6786
$ws.Range("J122").Value = 7996.6665
$ws.Range("K122").Value = 3800.0358
$ws.Range("L122").Value = 23989.9995
$ws.Range("M122").Value = -1350.0358
$ws.Range("N122").Value = -28889.9995
$ws.Range("H141").Value = 6265.625
$ws.Range("J141").Value = 4500
$ws.Range("L141").Value = 13500
$ws.Range("N141").Value = -23860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5195.297
$ws.Range("I32").Value = 2993.9856
$ws.Range("J32").Value = 12533
$ws.Range("K32").Value = 2993.9856
$ws.Range("L32").Value = 12533
$ws.Range("M32").Value = -2706.9856
$ws.Range("N32").Value = -13107
$ws.Range("H122").Value = 2241.96
$ws.Range("I122").Value = 1247.0714
$ws.Range("J122").Value = 3508.182
$ws.Range("K122").Value = 3741.2142
$ws.Range("L122").Value = 10524.546
$ws.Range("M122").Value = -1291.2142
$ws.Range("N122").Value = -15424.546
$ws.Range("H137").Value = 48510
$ws.Range("J137").Value = 48510
$ws.Range("L137").Value = 48510
$ws.Range("N137").Value = -58710

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1368.125
$ws.Range("I94").Value = 1442.1428
$ws.Range("J94").Value = 850
$ws.Range("K94").Value = 1442.1428
$ws.Range("L94").Value = 850
$ws.Range("M94").Value = -991.1428000000001
$ws.Range("N94").Value = -1752
$ws.Range("H107").Value = 1030.8235
$ws.Range("I107").Value = 948.1429000000001
$ws.Range("J107").Value = 1416.6666
$ws.Range("K107").Value = 948.1429000000001
$ws.Range("L107").Value = 1416.6666
$ws.Range("M107").Value = 971.8570999999999
$ws.Range("N107").Value = -5256.6666
$ws.Range("H137").Value = 32930
$ws.Range("J137").Value = 32930
$ws.Range("L137").Value = 32930
$ws.Range("N137").Value = -43130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 17762
$ws.Range("I36").Value = 15048
$ws.Range("J36").Value = 18666.666
$ws.Range("K36").Value = 15048
$ws.Range("L36").Value = 18666.666
$ws.Range("M36").Value = -14660
$ws.Range("N36").Value = -19442.666
$ws.Range("H40").Value = 17762
$ws.Range("I40").Value = 15048
$ws.Range("J40").Value = 18666.666
$ws.Range("K40").Value = 15048
$ws.Range("L40").Value = 18666.666
$ws.Range("M40").Value = -14888
$ws.Range("N40").Value = -18986.666
$ws.Range("H58").Value = 1751.7654
$ws.Range("I58").Value = 1547.8955
$ws.Range("J58").Value = 2727.4285
$ws.Range("K58").Value = 1547.8955
$ws.Range("L58").Value = 2727.4285
$ws.Range("M58").Value = -1344.8955
$ws.Range("N58").Value = -3133.4285
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H136").Value = 1751.7654
$ws.Range("I136").Value = 1547.8955
$ws.Range("J136").Value = 2727.4285
$ws.Range("K136").Value = 4643.6865
$ws.Range("L136").Value = 8182.2855
$ws.Range("M136").Value = -2093.6865
$ws.Range("N136").Value = -13282.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1653.9
$ws.Range("I41").Value = 493.33334
$ws.Range("J41").Value = 2151.2856
$ws.Range("K41").Value = 1480.00002
$ws.Range("L41").Value = 6453.8568
$ws.Range("M41").Value = -1142.00002
$ws.Range("N41").Value = -7129.8568
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864
$ws.Range("H107").Value = 23858144
$ws.Range("I107").Value = 295.84616
$ws.Range("J107").Value = 62627148
$ws.Range("K107").Value = 887.5384799999999
$ws.Range("L107").Value = 187881444
$ws.Range("M107").Value = 1032.46152
$ws.Range("N107").Value = -187885284
$ws.Range("H131").Value = 5955619
$ws.Range("J131").Value = 778.79486
$ws.Range("L131").Value = 2336.38458
$ws.Range("N131").Value = -12416.38458
$ws.Range("H137").Value = 3058.8
$ws.Range("I137").Value = 1390
$ws.Range("J137").Value = 4424.1816
$ws.Range("K137").Value = 4170
$ws.Range("L137").Value = 13272.5448
$ws.Range("M137").Value = 930
$ws.Range("N137").Value = -23472.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2440.3794
$ws.Range("I102").Value = 1971.55
$ws.Range("J102").Value = 3482.2222
$ws.Range("K102").Value = 1971.55
$ws.Range("L102").Value = 3482.2222
$ws.Range("M102").Value = -349.55
$ws.Range("N102").Value = -6726.2222
$ws.Range("H122").Value = 3679.4443
$ws.Range("I122").Value = 2443.8572
$ws.Range("J122").Value = 8004
$ws.Range("K122").Value = 7331.571599999999
$ws.Range("L122").Value = 24012
$ws.Range("M122").Value = -4881.571599999999
$ws.Range("N122").Value = -28912
$ws.Range("H123").Value = 10444.036
$ws.Range("J123").Value = 10444.036
$ws.Range("L123").Value = 10444.036
$ws.Range("N123").Value = -15344.036
$ws.Range("H124").Value = 43780
$ws.Range("J124").Value = 43780
$ws.Range("L124").Value = 43780
$ws.Range("N124").Value = -53600
$ws.Range("H132").Value = 2325.6924
$ws.Range("I132").Value = 1512.9688
$ws.Range("J132").Value = 6041
$ws.Range("K132").Value = 4538.9064
$ws.Range("L132").Value = 18123
$ws.Range("M132").Value = -2008.9064
$ws.Range("N132").Value = -23183
$ws.Range("H137").Value = 63601.6
$ws.Range("J137").Value = 63601.6
$ws.Range("L137").Value = 63601.6
$ws.Range("N137").Value = -73801.60000000001
$ws.Range("H141").Value = 54561.285
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 64385.8
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 64385.8
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -74745.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 14300
$ws.Range("J21").Value = 14300
$ws.Range("L21").Value = 14300
$ws.Range("N21").Value = -14648
$ws.Range("H40").Value = 9754.565000000001
$ws.Range("I40").Value = 9950.546
$ws.Range("J40").Value = 9574.916999999999
$ws.Range("K40").Value = 9950.546
$ws.Range("L40").Value = 9574.916999999999
$ws.Range("M40").Value = -9814.546
$ws.Range("N40").Value = -9846.916999999999
$ws.Range("H93").Value = 2114.6
$ws.Range("I93").Value = 1501.8889
$ws.Range("J93").Value = 3033.6667
$ws.Range("K93").Value = 1501.8889
$ws.Range("L93").Value = 3033.6667
$ws.Range("M93").Value = -253.8888999999999
$ws.Range("N93").Value = -5529.6667
$ws.Range("H122").Value = 5432.6665
$ws.Range("I122").Value = 2850
$ws.Range("K122").Value = 8550
$ws.Range("M122").Value = -6100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 334.72
$ws.Range("I113").Value = 281.55554
$ws.Range("J113").Value = 471.42856
$ws.Range("K113").Value = 844.66662
$ws.Range("L113").Value = 1414.28568
$ws.Range("M113").Value = 1325.33338
$ws.Range("N113").Value = -5754.28568
$ws.Range("H122").Value = 3442.4688
$ws.Range("I122").Value = 2255.1904
$ws.Range("K122").Value = 6765.5712
$ws.Range("M122").Value = -4315.5712
$ws.Range("H140").Value = 53721.5
$ws.Range("J140").Value = 53721.5
$ws.Range("L140").Value = 53721.5
$ws.Range("N140").Value = -64081.5
$ws.Range("H141").Value = 43635.355
$ws.Range("J141").Value = 43635.355
$ws.Range("L141").Value = 43635.355
$ws.Range("N141").Value = -53995.355
